$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '37.371.41'
$ws.Range('E2').Value = '  +0.54%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.016.97'
$ws.Range('E3').Value = '  +0.74%  '

$ws.Range('E4').Value = '  -0.12%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '259.81'
$ws.Range('E5').Value = '  +5.67%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.610'
$ws.Range('E6').Value = '  -2.13%  '

$ws.Range('E7').Value = '  -0.04%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '56.50'
$ws.Range('E8').Value = '  -5.82%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.388'
$ws.Range('E9').Value = '  +0.72%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0774'
$ws.Range('E10').Value = '  -3.82%  '

$ws.Range('E11').Value = '  -1.87%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '14.34'
$ws.Range('E12').Value = '  -4.27%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.310.29'
$ws.Range('E13').Value = '  +0.62%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.808'
$ws.Range('E14').Value = '  -3.94%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '21.04'
$ws.Range('E15').Value = '  -6.55%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.28'
$ws.Range('E16').Value = '  -2.81%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.022.38'
$ws.Range('E17').Value = '  +1.42%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '37.320.68'
$ws.Range('E18').Value = '  +0.55%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '69.92'
$ws.Range('E19').Value = '  -0.34%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0839'
$ws.Range('E20').Value = '  -2.63%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.20'
$ws.Range('E21').Value = '  +0.51%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '228.66'
$ws.Range('E22').Value = '  -0.63%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.66'
$ws.Range('E23').Value = '  +7.72%  '

$ws.Range('E24').Value = '  +0.03%  '

$ws.Range('E25').Value = '  +0.62%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '165.02'
$ws.Range('E26').Value = '  +0.66%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.00'
$ws.Range('E27').Value = '  -4.70%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.74'
$ws.Range('E28').Value = '  +0.57%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.131'
$ws.Range('E29').Value = '  -7.27%  '

$ws.Range('E30').Value = '  -2.35%  '

$ws.Range('E31').Value = '  -0.60%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.68'
$ws.Range('E32').Value = '  -2.47%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0651'
$ws.Range('E33').Value = '  -1.06%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.64'
$ws.Range('E34').Value = '  +3.43%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.41'
$ws.Range('E35').Value = '  +1.60%  '

$ws.Range('E36').Value = '  +1.11%  '

$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.39'
$ws.Range('E37').Value = '  +2.70%  '

$ws.Range('B38').Value = 'BinanceUSD'
$ws.Range('C38').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.00'
$ws.Range('E38').Value = '  -0.06%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.26'
$ws.Range('E39').Value = '  -2.22%  '

$ws.Range('E40').Value = '  +3.94%  '

$ws.Range('E41').Value = '  +2.59%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0214'
$ws.Range('E42').Value = '  -0.28%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0934'
$ws.Range('E43').Value = '  -4.96%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.400.76'
$ws.Range('E44').Value = '  +2.34%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '90.35'
$ws.Range('E45').Value = '  -0.89%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '15.82'
$ws.Range('E46').Value = '  -4.66%  '

$ws.Range('E47').Value = '  -1.54%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.13'
$ws.Range('E48').Value = '  -3.32%  '

$ws.Range('E49').Value = '  +1.95%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.201.94'
$ws.Range('E50').Value = '  +0.62%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.96'
$ws.Range('E51').Value = '  -5.51%  '
